$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "Precio del productooooooo"
$ws.Range("K2").Value = "2023/05/19"
$ws.Range("K3").Value = " 15:18"
